# Updates betting odds values in the Betfair Back/Lay sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$changes = @{
    "N2"   = 3.85
    "P2"   = 1.99
    "Q2"   = 1.95
    "R2"   = 1.38
    "U2"   = 2.14

    "F3"   = 2.58
    "G3"   = 2.74
    "H3"   = 2.6
    "I3"   = 2.74
    "L3"   = 1.29
    "P3"   = 2.6
    "R3"   = 1.66
    "S3"   = 2.4
    "T3"   = 1.53
    "U3"   = 2.72
    "V3"   = 1.57
    "W3"   = 1.58
    "X3"   = 26
    "Y3"   = 17
    "Z3"   = 22
    "AB3"  = 17
    "AD3"  = 13
    "AF3"  = 22
    "AK3"  = 980
    "AL3"  = 34
    "AN3"  = 14.5
    "AO3"  = 14.5

    "H5"   = 1.91
    "I5"   = 2.38
    "K5"   = 3.8
    "P5"   = 1.35
    "Q5"   = 1.01

    "F6"   = 1.67

    "F7"   = 6.6
    "G7"   = 1000
    "I7"   = 1.31
    "J7"   = 5.2
    "K7"   = 950

    "F8"   = 1.93
    "G8"   = 2.42
    "J8"   = 2.58
    "K8"   = 3.8

    "N11"  = 4.9

    "P12"  = 1.53
    "R12"  = 1.19
    "AM12" = 270
    "AN12" = 30

    "N13"  = 2.74
    "P13"  = 1.58
    "U13"  = 1.76
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
